$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 48; this shifts existing rows 48-137 down to 49-138.
$ws.Rows.Item(48).Insert()

# Populate the newly inserted row 48 with the new data record.
$ws.Range("A48").Value = 10
$ws.Range("B48").Value = "Vega Modelo de Temuco"
$ws.Range("C48").Value = "La Araucanía"
$ws.Range("D48").Value = 45281
$ws.Range("D48").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E48").Value = 9
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100108
$ws.Range("H48").Value = "Tropicales y subtropicales"
$ws.Range("I48").Value = 100108007
$ws.Range("J48").Value = "Coco"
$ws.Range("K48").Value = "Sin especificar"
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 120
$ws.Range("N48").Value = 30000
$ws.Range("O48").Value = 32000
$ws.Range("P48").Value = 30833
$ws.Range("Q48").Value = "$/malla 20 unidades"
$ws.Range("R48").Value = "Perú"
$ws.Range("S48").Value = 1542
$ws.Range("T48").Value = 20
